# Refresh the cryptos price/volume snapshot (GitHub Actions style update).
# Column D = Price (text, since some values use "."-grouped thousands),
# Column E = Volume(1h) change, stored as "  +X.XX%  " padded text.
# A leading "'" forces cells that look like plain decimals (e.g. "239.12")
# to stay text instead of being auto-parsed as numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.470.14"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.872.28"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'0.7179"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Value = "'239.12"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07822"
$ws.Range("E8").Value = "  -4.41%  "
$ws.Range("D9").Value = "'0.3068"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'25.27"
$ws.Range("E10").Value = "  +8.81%  "
$ws.Range("D11").Value = "'0.08235"
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.880.83"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "'5.228"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'0.7205"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "'89.93"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "29.520.06"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'0.000007850"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "'240.31"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "'13.28"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "2.135.81"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'7.740"
$ws.Range("E24").Value = "  +4.14%  "
$ws.Range("D25").Value = "'0.1552"
$ws.Range("E25").Value = "  +5.80%  "
$ws.Range("D26").Value = "'162.65"
$ws.Range("D27").Value = "'8.974"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'18.29"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "'1.930"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").Value = "'1.357"
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("D31").Value = "'1.482"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'4.327"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").Value = "'4.076"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "'0.05244"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "'1.197"
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("D36").Value = "'0.7164"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "'0.9995"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").Value = "'2.676"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'0.01866"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "'2.718"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").Value = "1.175.64"
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").Value = "'0.9061"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").Value = "'5.986"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").Value = "'0.4305"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "'71.38"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'102.21"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "'0.5367"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  -0.79%  "
$ws.Range("D50").Value = "'9.146"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").Value = "'7.010"
$ws.Range("E51").Value = "  +0.73%  "